$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.725.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '''3.099.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.67%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").Value = '''542.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.73%  '

$ws.Range("D6").Value = '''137.10'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.85%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '''3.092.91'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.81%  '

$ws.Range("D9").Value = '''0.498'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.41%  '

$ws.Range("D10").Value = '''0.157'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.22%  '

$ws.Range("D11").Value = '''6.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.45%  '

$ws.Range("D12").Value = '''0.461'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.43%  '

$ws.Range("D13").Value = '''0.0000228'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.94%  '

$ws.Range("D14").Value = '''35.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.89%  '

$ws.Range("D15").Value = '''3.592.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.31%  '

$ws.Range("D16").Value = '''63.557.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").Value = '''3.089.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("D19").Value = '''6.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("D20").Value = '''491.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.60%  '

$ws.Range("D21").Value = '''13.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.23%  '

$ws.Range("D22").Value = '''0.706'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.02%  '

$ws.Range("D23").Value = '''7.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.73%  '

$ws.Range("D24").Value = '''79.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.85%  '

$ws.Range("D25").Value = '''12.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.92%  '

$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '''8.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.50%  '

$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '''2.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.08%  '

$ws.Range("D29").Value = '''0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.44%  '

$ws.Range("D30").Value = '''26.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.41%  '

$ws.Range("D31").Value = '''1.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.78%  '

$ws.Range("E32").Value = '  +1.22%  '

$ws.Range("D33").Value = '''2.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.31%  '

$ws.Range("D34").Value = '''57.30'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.12%  '

$ws.Range("D35").Value = '''5.43'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.98%  '

$ws.Range("D36").Value = '''6.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.04%  '

$ws.Range("D37").Value = '''494.93'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.69%  '

$ws.Range("D38").Value = '''3.192.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.10%  '

$ws.Range("D39").Value = '''0.0403'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.17%  '

$ws.Range("D40").Value = '''0.0807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.01%  '

$ws.Range("D41").Value = '''0.118'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.51%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '''2.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.54%  '

$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Value = '''8.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.60%  '

$ws.Range("D44").Value = '''0.259'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.09%  '

$ws.Range("D46").Value = '''0.0₃0549'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.47%  '

$ws.Range("D47").Value = '''2.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.31%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '''24.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.33%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''121.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("E50").Value = '  +3.39%  '

$ws.Range("D51").Value = '''2.38'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.47%  '
